$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 235.28
$ws.Range("I33").Value = 125.625
$ws.Range("J33").Value = 430.22223
$ws.Range("K33").Value = 125.625
$ws.Range("L33").Value = 430.22223
$ws.Range("M33").Value = 103.375
$ws.Range("N33").Value = -888.2222300000001

$ws.Range("H40").Value = 2531.1765
$ws.Range("I40").Value = 2050
$ws.Range("J40").Value = 2868
$ws.Range("K40").Value = 2050
$ws.Range("L40").Value = 2868
$ws.Range("M40").Value = -1875
$ws.Range("N40").Value = -3218

$ws.Range("H98").Value = 1303.6923
$ws.Range("I98").Value = 539.9
$ws.Range("J98").Value = 3849.6667
$ws.Range("K98").Value = 539.9
$ws.Range("L98").Value = 3849.6667
$ws.Range("M98").Value = 958.1

$ws.Range("H103").Value = 13334200
$ws.Range("I103").Value = 699.8333
$ws.Range("J103").Value = 22223200
$ws.Range("K103").Value = 2099.4999
$ws.Range("L103").Value = 66669600
$ws.Range("M103").Value = -1513.4999
$ws.Range("N103").Value = -66670772

$ws.Range("H122").Value = 1303.6923
$ws.Range("I122").Value = 539.9
$ws.Range("J122").Value = 3849.6667
$ws.Range("K122").Value = 1619.7
$ws.Range("L122").Value = 11549.0001
$ws.Range("M122").Value = 830.3000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1932.45
$ws.Range("I97").Value = 1153.1072
$ws.Range("J97").Value = 3750.9167
$ws.Range("K97").Value = 1153.1072
$ws.Range("L97").Value = 3750.9167
$ws.Range("M97").Value = -657.1071999999999
$ws.Range("N97").Value = -4742.9167

$ws.Range("H102").Value = 166668740
$ws.Range("I102").Value = 2478
$ws.Range("J102").Value = 1000000000
$ws.Range("K102").Value = 2478
$ws.Range("L102").Value = 1000000000
$ws.Range("M102").Value = -856
$ws.Range("N102").Value = -1000003244

$ws.Range("H132").Value = 203242.48
$ws.Range("I132").Value = 264488.1
$ws.Range("J132").Value = 9298
$ws.Range("K132").Value = 793464.2999999999
$ws.Range("L132").Value = 27894
$ws.Range("M132").Value = -790934.2999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9262517
$ws.Range("I20").Value = 15877052
$ws.Range("J20").Value = 2167.8667
$ws.Range("K20").Value = 15877052
$ws.Range("L20").Value = 2167.8667
$ws.Range("M20").Value = -15876805

$ws.Range("H94").Value = 1149.0526
$ws.Range("I94").Value = 1071.6923
$ws.Range("J94").Value = 1316.6666
$ws.Range("K94").Value = 1071.6923
$ws.Range("L94").Value = 1316.6666
$ws.Range("M94").Value = -620.6922999999999
$ws.Range("N94").Value = -2218.6666

$ws.Range("H105").Value = 2985.875
$ws.Range("I105").Value = 2189.7222
$ws.Range("J105").Value = 4009.5
$ws.Range("K105").Value = 2189.7222
$ws.Range("L105").Value = 4009.5
$ws.Range("M105").Value = -442.7222000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2174.5527
$ws.Range("I31").Value = 1624.5454
$ws.Range("J31").Value = 2930.8125
$ws.Range("K31").Value = 1624.5454
$ws.Range("L31").Value = 2930.8125
$ws.Range("M31").Value = -1329.5454
$ws.Range("N31").Value = -3520.8125

$ws.Range("H34").Value = 2174.5527
$ws.Range("I34").Value = 1624.5454
$ws.Range("J34").Value = 2930.8125
$ws.Range("K34").Value = 1624.5454
$ws.Range("L34").Value = 2930.8125
$ws.Range("M34").Value = -1422.5454
$ws.Range("N34").Value = -3334.8125

$ws.Range("H58").Value = 1291.1666
$ws.Range("I58").Value = 1203.1852
$ws.Range("J58").Value = 2083
$ws.Range("K58").Value = 1203.1852
$ws.Range("L58").Value = 2083
$ws.Range("M58").Value = -1000.1852
$ws.Range("N58").Value = -2489

$ws.Range("H62").Value = 2534.5386
$ws.Range("I62").Value = 2495.3635
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 2495.3635
$ws.Range("L62").Value = 2750
$ws.Range("M62").Value = -1871.3635
$ws.Range("N62").Value = -3998

$ws.Range("H65").Value = 2534.5386
$ws.Range("I65").Value = 2495.3635
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 12476.8175
$ws.Range("L65").Value = 13750
$ws.Range("M65").Value = -9356.817499999999
$ws.Range("N65").Value = -19990

$ws.Range("H132").Value = 2525.9524
$ws.Range("I132").Value = 2034.1875
$ws.Range("J132").Value = 4099.6
$ws.Range("K132").Value = 6102.5625
$ws.Range("L132").Value = 12298.8
$ws.Range("M132").Value = -3572.5625

$ws.Range("H136").Value = 1291.1666
$ws.Range("I136").Value = 1203.1852
$ws.Range("J136").Value = 2083
$ws.Range("K136").Value = 3609.5556
$ws.Range("L136").Value = 6249
$ws.Range("M136").Value = -1059.5556
$ws.Range("N136").Value = -11349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H113").Value = 567.4
$ws.Range("I113").Value = 518.5789
$ws.Range("J113").Value = 611.5714
$ws.Range("K113").Value = 1555.7367
$ws.Range("L113").Value = 1834.7142
$ws.Range("M113").Value = 614.2633000000001
$ws.Range("N113").Value = -6174.7142

$ws.Range("H121").Value = 534253.4399999999
$ws.Range("I121").Value = 376.66666
$ws.Range("J121").Value = 607054.8
$ws.Range("K121").Value = 1129.99998
$ws.Range("L121").Value = 1821164.4
$ws.Range("M121").Value = 180.0000199999999
$ws.Range("N121").Value = -1823784.4

$ws.Range("H122").Value = 2844.54
$ws.Range("I122").Value = 213.8125
$ws.Range("J122").Value = 3345.6309
$ws.Range("K122").Value = 1924.3125
$ws.Range("L122").Value = 30110.6781
$ws.Range("M122").Value = 525.6875
$ws.Range("N122").Value = -35010.6781

$ws.Range("H129").Value = 1884.7059
$ws.Range("I129").Value = 960
$ws.Range("J129").Value = 2706.6667
$ws.Range("K129").Value = 2880
$ws.Range("L129").Value = 8120.000100000001
$ws.Range("M129").Value = 2120
$ws.Range("N129").Value = -18120.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.90909000000001
$ws.Range("I2").Value = 43.333332
$ws.Range("J2").Value = 74.375
$ws.Range("K2").Value = 43.333332
$ws.Range("L2").Value = 74.375
$ws.Range("M2").Value = 69.666668
$ws.Range("N2").Value = -300.375

$ws.Range("H97").Value = 1643.3334
$ws.Range("I97").Value = 1643.3334
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1643.3334
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1147.3334
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 880.625
$ws.Range("I16").Value = 883.2857
$ws.Range("J16").Value = 862
$ws.Range("K16").Value = 883.2857
$ws.Range("L16").Value = 862
$ws.Range("M16").Value = -713.2857
$ws.Range("N16").Value = -1202

$ws.Range("H100").Value = 2278.2354
$ws.Range("I100").Value = 2493
$ws.Range("J100").Value = 1971.4286
$ws.Range("K100").Value = 2493
$ws.Range("L100").Value = 1971.4286
$ws.Range("M100").Value = -1952

$ws.Range("H133").Value = 60797.6
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 60797.6
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 60797.6
$ws.Range("N133").Value = -65857.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1144.1923
$ws.Range("I126").Value = 739.5417
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 2218.6251
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = 251.3748999999998
